# Updates the cryptocurrency price/volume table (rows 2-51) on Sheet1 to the
# latest scraped values. Column D ("Price") holds display-formatted price
# strings (thousand separators as literal dots, trailing zeros, subscript
# notation, etc.) that must remain plain TEXT rather than being
# reinterpreted as numbers -- so those cells are written via a
# leading-apostrophe Formula assignment (the same thing Excel does when a
# user types '123.45 into a cell), which keeps the cell type Text and
# preserves the exact digit string. Column E ("Volume(1h)") values are
# padded percentage strings and are assigned directly via .Value (they
# already aren't numeric due to the surrounding spaces / % sign).
# Rows 31 and 32 also swap which coin (Mantle / EthereumClassic) occupies
# which rank, so their Coin name (B) and Link (C) columns are rewritten too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").Formula = "'61.177.88"
$ws.Range("E2").Value = '  -2.49%  '

# Row 3
$ws.Range("D3").Formula = "'3.004.77"
$ws.Range("E3").Value = '  -1.82%  '

# Row 4
$ws.Range("D4").Formula = "'0.999"
$ws.Range("E4").Value = '  -0.12%  '

# Row 5
$ws.Range("D5").Formula = "'535.36"
$ws.Range("E5").Value = '  -0.01%  '

# Row 6
$ws.Range("D6").Formula = "'134.76"
$ws.Range("E6").Value = '  +1.84%  '

# Row 7
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").Formula = "'2.999.15"
$ws.Range("E8").Value = '  -1.73%  '

# Row 9
$ws.Range("D9").Formula = "'0.499"
$ws.Range("E9").Value = '  +1.46%  '

# Row 10
$ws.Range("D10").Formula = "'0.149"
$ws.Range("E10").Value = '  -3.02%  '

# Row 11
$ws.Range("D11").Formula = "'6.12"
$ws.Range("E11").Value = '  +1.07%  '

# Row 12
$ws.Range("D12").Formula = "'0.447"
$ws.Range("E12").Value = '  -0.41%  '

# Row 13
$ws.Range("D13").Formula = "'0.0000221"
$ws.Range("E13").Value = '  -0.64%  '

# Row 14
$ws.Range("D14").Formula = "'34.34"
$ws.Range("E14").Value = '  +0.94%  '

# Row 15
$ws.Range("D15").Formula = "'3.487.53"
$ws.Range("E15").Value = '  -1.83%  '

# Row 16
$ws.Range("E16").Value = '  -0.22%  '

# Row 17
$ws.Range("D17").Formula = "'61.153.04"
$ws.Range("E17").Value = '  -2.53%  '

# Row 18
$ws.Range("D18").Formula = "'2.999.48"
$ws.Range("E18").Value = '  -2.02%  '

# Row 19
$ws.Range("D19").Formula = "'6.63"
$ws.Range("E19").Value = '  +0.59%  '

# Row 20
$ws.Range("D20").Formula = "'465.42"
$ws.Range("E20").Value = '  -3.02%  '

# Row 21
$ws.Range("D21").Formula = "'13.25"
$ws.Range("E21").Value = '  +0.02%  '

# Row 22
$ws.Range("D22").Formula = "'0.677"
$ws.Range("E22").Value = '  -1.83%  '

# Row 23
$ws.Range("D23").Formula = "'6.98"
$ws.Range("E23").Value = '  -1.13%  '

# Row 24
$ws.Range("D24").Formula = "'79.54"
$ws.Range("E24").Value = '  +1.06%  '

# Row 25
$ws.Range("D25").Formula = "'12.10"
$ws.Range("E25").Value = '  +0.83%  '

# Row 26
$ws.Range("E26").Value = '  +0.04%  '

# Row 27
$ws.Range("D27").Formula = "'2.68"
$ws.Range("E27").Value = '  -0.11%  '

# Row 28
$ws.Range("D28").Formula = "'7.91"
$ws.Range("E28").Value = '  -1.37%  '

# Row 29
$ws.Range("D29").Formula = "'0.997"
$ws.Range("E29").Value = '  -0.07%  '

# Row 30
$ws.Range("D30").Formula = "'1.91"
$ws.Range("E30").Value = '  +2.42%  '

# Row 31
$ws.Range("B31").Value = 'Mantle'
$ws.Range("C31").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D31").Formula = "'1.15"
$ws.Range("E31").Value = '  +4.31%  '

# Row 32
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Formula = "'25.59"
$ws.Range("E32").Value = '  -1.08%  '

# Row 33
$ws.Range("D33").Formula = "'5.52"
$ws.Range("E33").Value = '  +3.64%  '

# Row 34
$ws.Range("D34").Formula = "'55.64"
$ws.Range("E34").Value = '  -2.12%  '

# Row 35
$ws.Range("D35").Formula = "'2.29"
$ws.Range("E35").Value = '  -2.05%  '

# Row 36
$ws.Range("D36").Formula = "'5.91"
$ws.Range("E36").Value = '  -0.89%  '

# Row 37
$ws.Range("D37").Formula = "'460.36"
$ws.Range("E37").Value = '  -3.15%  '

# Row 38
$ws.Range("D38").Formula = "'3.204.36"
$ws.Range("E38").Value = '  +3.91%  '

# Row 39
$ws.Range("D39").Formula = "'0.0790"
$ws.Range("E39").Value = '  +0.10%  '

# Row 40
$ws.Range("D40").Formula = "'0.0385"
$ws.Range("E40").Value = '  -1.55%  '

# Row 41
$ws.Range("E41").Value = '  +2.60%  '

# Row 42
$ws.Range("D42").Formula = "'8.17"
$ws.Range("E42").Value = '  +1.66%  '

# Row 43
$ws.Range("D43").Formula = "'27.82"
$ws.Range("E43").Value = '  +14.83%  '

# Row 44
$ws.Range("D44").Formula = "'2.48"
$ws.Range("E44").Value = '  -4.49%  '

# Row 46
$ws.Range("D46").Formula = "'0.247"
$ws.Range("E46").Value = '  -1.29%  '

# Row 47
$ws.Range("D47").Formula = "'2.01"
$ws.Range("E47").Value = '  +0.77%  '

# Row 48
$ws.Range("D48").Formula = "'119.34"
$ws.Range("E48").Value = '  -1.35%  '

# Row 49
$ws.Range("D49").Formula = "'0.109"
$ws.Range("E49").Value = '  +0.96%  '

# Row 50
$ws.Range("D50").Formula = "'0.0₃0495"
$ws.Range("E50").Value = '  -7.07%  '

# Row 51
$ws.Range("E51").Value = '  +8.12%  '

Write-Output "Updated cryptos list"
